$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the environment-specific values from "test18" to "test7"
$ws.Range("A2").Value = "https://test7.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test7.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test7.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest7"
$ws.Range("G2").Value = "test7"
$ws.Range("K2").Value = "test7"

# Update the active selection shown in the sheet view
$ws.Range("C18").Select()
